$wb = $excel.ActiveWorkbook

# --- Sheet "General" (sheet1.xml) ---------------------------------------
$general = $wb.Worksheets.Item("General")
$general.Activate()
$general.Range("D3").Select()

# --- Sheet "Gana" (sheet2.xml) ------------------------------------------
$gana = $wb.Worksheets.Item("Gana")
$gana.Activate()

# BDD step text rewrites: drop the quotes around the literal value and
# fix the "enter"/"Click" casing (now driven from an enum, per the
# commit message), then clear the now-unused XPath lookup cells.
# Values are written in this order so newly-created shared-string
# entries land in the same order the original workbook produced them.
$gana.Range("C5").Value = 'Enter reddy into "Your Last Name"'
$gana.Range("C6").Value = 'Enter 08999999 into "Phone Number"'
$gana.Range("C7").Value = 'Enter dublin into "Your Address"'
$gana.Range("C8").Value = 'Enter ireland into "Your Address Two"'
$gana.Range("C9").Value = 'Enter 12344 into "Your Pin Code"'
$gana.Range("C3").Value = 'Enter madan into "Your Name"'
$gana.Range("C4").Value = 'Enter mohan into "Your Middle Name"'
$gana.Range("C13").Value = 'click on "Never Registered"'

$gana.Range("D3").ClearContents()
$gana.Range("D4").ClearContents()
$gana.Range("D5").ClearContents()
$gana.Range("D6").ClearContents()
$gana.Range("D7").ClearContents()
$gana.Range("D8").ClearContents()
$gana.Range("D9").ClearContents()
$gana.Range("D10").ClearContents()
$gana.Range("D11").ClearContents()
$gana.Range("D13").ClearContents()

# Row 1 no longer needs the taller custom height - restore to the sheet
# default.
$gana.Rows.Item(1).AutoFit()

# New column F gets an explicit width.
$gana.Range("F1").ColumnWidth = 46.14

# Selection moves to D13 on this sheet.
$gana.Range("D13").Select()
